$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Faith Pleases God" (10/9/2025) entry spanning rows 126-128 is being
# replaced with a "Guest Speaker Week (Unknown Topic)" entry dated 11/2/2025.
# Row 126 also swaps which song/CCLI pairing it carries with row 127, and the
# previously-blank CCLI cell on that row now gets a value.
# ---------------------------------------------------------------------------

# --- Row 126: now carries the song/CCLI that used to be on row 127, but with
#     the "accent" (red) formatting that the other row used to have -------
$ws.Range("B121").Copy()
$ws.Range("B126").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A126").Value = "11/2/2025"
$ws.Range("B126").Value = "#27 - Kuv Twb Raug Kev Txomnyem Ntsuav"
$ws.Range("C126").Value = "No record to report"
$ws.Range("D126").Value = "HBNA Songbook"
$ws.Range("E126").Value = "Guest Speaker Week (Unknown Topic)"

# --- Row 127: now carries the song/CCLI that used to be on row 126, with the
#     "highlight" (orange) formatting used elsewhere for that song --------
$ws.Range("B119").Copy()
$ws.Range("B127").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A127").Value = "11/2/2025"
$ws.Range("B127").Value = "#23 - Kuv Muab Siab Npuab Yexus"
$ws.Range("C127").Value = "Public - No reporting needed"
$ws.Range("D127").Value = "HBNA Songbook"
$ws.Range("E127").Value = "Guest Speaker Week (Unknown Topic)"

# --- Row 128: same song/CCLI as before, just the new date + topic --------
$ws.Range("A128").Value = "11/2/2025"
$ws.Range("B128").Value = "#57 - Kuv Tus Kwvluag Yog Tswv Yexus"
$ws.Range("C128").Value = "Public - No reporting needed"
$ws.Range("D128").Value = "HBNA Songbook"
$ws.Range("E128").Value = "Guest Speaker Week (Unknown Topic)"

# --- Sheet view: scroll the frozen pane down and move the active selection
#     to reflect where the author left off editing ------------------------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 40
$ws.Range("B132").Select()
